$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values stay as text (matches source data which stores
# prices/percentages as text strings, not numbers)
$textCells = @("D2","E2","D3","E3","D4","E4","D5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D26","E26","E27","E28","E40","D41","E41","E42","E43","D44","E44","E45","E46","D47","E47","D48","E48","E49","E50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update coin name / link cells (plain text, shifted rows)
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B18").Value = "UpBots"
$ws.Range("C18").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# Update price / volume cells (kept as text)
$ws.Range("D2").Value = "246.47"
$ws.Range("E2").Value = "0.57%"
$ws.Range("D3").Value = "29.71"
$ws.Range("E3").Value = "10.26%"
$ws.Range("D4").Value = "5.169"
$ws.Range("E4").Value = "1.96%"
$ws.Range("D5").Value = "0.05705"
$ws.Range("D6").Value = "6.613"
$ws.Range("E6").Value = "2.12%"
$ws.Range("D7").Value = "3.072"
$ws.Range("E7").Value = "2.34%"
$ws.Range("D8").Value = "0.8579"
$ws.Range("E8").Value = "4.61%"
$ws.Range("D9").Value = "0.8679"
$ws.Range("E9").Value = "3.29%"
$ws.Range("D10").Value = "0.1364"
$ws.Range("E10").Value = "2.63%"
$ws.Range("D11").Value = "0.07099"
$ws.Range("E11").Value = "2.70%"
$ws.Range("D12").Value = "0.02921"
$ws.Range("E12").Value = "2.14%"
$ws.Range("D13").Value = "0.09382"
$ws.Range("E13").Value = "-0.14%"
$ws.Range("D14").Value = "0.001509"
$ws.Range("E14").Value = "-0.06%"
$ws.Range("D15").Value = "0.04170"
$ws.Range("E15").Value = "1.94%"
$ws.Range("D16").Value = "0.0005991"
$ws.Range("E16").Value = "0.08%"
$ws.Range("D17").Value = "0.006085"
$ws.Range("E17").Value = "0.94%"
$ws.Range("D18").Value = "0.007491"
$ws.Range("E18").Value = "0.07%"
$ws.Range("D19").Value = "3.490"
$ws.Range("E19").Value = "-0.52%"
$ws.Range("D20").Value = "2.275"
$ws.Range("E20").Value = "2.14%"
$ws.Range("D21").Value = "0.3175"
$ws.Range("E21").Value = "-0.07%"
$ws.Range("D22").Value = "0.03323"
$ws.Range("E22").Value = "4.68%"
$ws.Range("D23").Value = "0.1325"
$ws.Range("E23").Value = "2.13%"
$ws.Range("D24").Value = "3.475"
$ws.Range("E24").Value = "-2.36%"
$ws.Range("D26").Value = "0.005031"
$ws.Range("E26").Value = "26.92%"
$ws.Range("E27").Value = "0.42%"
$ws.Range("E28").Value = "23.60%"
$ws.Range("E40").Value = "1.49%"
$ws.Range("D41").Value = "0.005771"
$ws.Range("E41").Value = "68.01%"
$ws.Range("E42").Value = "1.33%"
$ws.Range("E43").Value = "-14.27%"
$ws.Range("D44").Value = "0.009977"
$ws.Range("E44").Value = "6.38%"
$ws.Range("E45").Value = "-0.05%"
$ws.Range("E46").Value = "0.09%"
$ws.Range("D47").Value = "0.06001"
$ws.Range("E47").Value = "-40.83%"
$ws.Range("D48").Value = "0.002564"
$ws.Range("E48").Value = "-1.11%"
$ws.Range("E49").Value = "0.09%"
$ws.Range("E50").Value = "0.09%"
